$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
Write-Host $ws.Range("A1").Value
Write-Host $ws.Range("D1").Value
